$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 5 and 6 have their D, M, N, O, P, Q, R, S, T values swapped.

$cols = @("D", "M", "N", "O", "P", "Q", "R", "S", "T")

foreach ($col in $cols) {
    $addr5 = $col + "5"
    $addr6 = $col + "6"
    $val5 = $ws.Range($addr5).Value2
    $val6 = $ws.Range($addr6).Value2
    $ws.Range($addr5).Value2 = $val6
    $ws.Range($addr6).Value2 = $val5
}
